$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The progress table on the slide ("표 4") is shape #2.
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table

# Week 8 row ("사운드 / 게임종료", "밸런스 조절 / 최적화"): bump completion 50% -> 90%
$soundCell = $tbl.Cell(8, 4)
$soundCell.Shape.TextFrame.TextRange.Text = "90%"

# Week 9 row ("추가범위 구현"): record a completion value of 60% (cell was empty)
$extraCell = $tbl.Cell(9, 4)
$extraCell.Shape.TextFrame.TextRange.Text = "60%"
